$d = $word.ActiveDocument

# --- Change 1: merge the two "SAT Dec 9" / " 10:14:58 PST 2017" runs into one run ---
$rng1 = $d.Content
$rng1.Find.Execute("SAT Dec 9 10:14:58 PST 2017", $false, $false, $false, $false, $false, $true, 1, $false, "SAT Dec 9 10:14:58 PST 2017", 2) | Out-Null

# --- Change 2: insert a new purchase-entry block right after the last "- CASH" paragraph ---
$rng2 = $d.Content
$count = 0
while ($rng2.Find.Execute("- CASH", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $count++
    $rng2.Collapse(0)
    if ($count -ge 13) { break }
}

$p = $rng2.Paragraphs(1)
$insertPos = $p.Range.End
$target = $d.Range($insertPos, $insertPos)

$newBlockXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>THU Dec 14</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t xml:space=`"preserve`"> 10:49:19 PST 2017</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Person Name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- TRH</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Bill number</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 1948</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>---------------------------------------------------------------</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Item Name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- CARROT</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Number of Pockets</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 1</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Number of KGs</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 35</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Rate</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 42</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Transport &amp; Miscellaneous</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 10</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Total Price</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 1480.0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr><w:t>Amount balance</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr><w:tab/><w:t>- 5744.0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr></w:p><w:p></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$target.InsertXML($newBlockXml)

Write-Output "Change1 matches merged; Change2 inserted after match #$count (pos $insertPos)"
